$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rewrite the stock table -------------------------------------------
# Original data occupied A1:C7 (with several formulas in column B).
# The new data occupies A2:C11, all as plain values (no formulas).

# Clear old row 1 entirely so it disappears from the sheet.
$ws.Range("A1:C1").ClearContents()

# Row 2
$ws.Range("A2").Value = 940
$ws.Range("B2").Value = 2400
$ws.Range("C2").Value = 70598

# Row 3
$ws.Range("A3").Value = 941
$ws.Range("B3").Value = 1920
$ws.Range("C3").Value = "MANAOS"

# Row 4
$ws.Range("A4").Value = 942
$ws.Range("B4").Value = 1920
$ws.Range("C4").Value = "IM"

# Row 5 (also clear any leftover C5 from the old layout)
$ws.Range("A5").Value = 944
$ws.Range("B5").Value = 1920
$ws.Range("C5").ClearContents()

# Row 6
$ws.Range("A6").Value = 970
$ws.Range("B6").Value = 168

# Row 7
$ws.Range("A7").Value = 965
$ws.Range("B7").Value = 960

# Row 8
$ws.Range("A8").Value = 9402
$ws.Range("B8").Value = 480

# Row 9
$ws.Range("A9").Value = 946
$ws.Range("B9").Value = 480

# Row 10
$ws.Range("A10").Value = 950
$ws.Range("B10").Value = 480

# Row 11
$ws.Range("A11").Value = 920
$ws.Range("B11").Value = 450

# --- Selection -----------------------------------------------------------
$ws.Range("C5").Select() | Out-Null

# --- AutoFilter + hidden _FilterDatabase defined name ---------------------
$ws.Range("A1:S1").AutoFilter() | Out-Null
$fd = $ws.Names.Add("_xlnm._FilterDatabase", "=Sheet1!`$B`$2:`$B`$2")
$fd.RefersTo = "=Sheet1!`$A`$1:`$S`$1"
$fd.Visible = $false
